$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Containers")

# Update the description of the "wallpapers" container (row 7, column C)
$ws.Range("C7").Value = "All of the wallpaper images that could appear as the background of the home screen. This folder still exists but is not being used anymore."

# Update the active selection to C8 (as in the target diff)
$ws.Range("C8").Select()
